$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.844.60"
Set-TextValue "E2" "  -1.07%  "
Set-TextValue "D3" "1.733.01"
Set-TextValue "E3" "  +0.83%  "
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  +0.09%  "
Set-TextValue "D5" "310.75"
Set-TextValue "E5" "  -0.31%  "
Set-TextValue "E6" "  -0.16%  "
Set-TextValue "D7" "0.5010"
Set-TextValue "E7" "  +8.38%  "
Set-TextValue "D8" "0.3569"
Set-TextValue "E8" "  +4.36%  "
Set-TextValue "D9" "42.24"
Set-TextValue "E9" "  +0.83%  "
Set-TextValue "D10" "0.07235"
Set-TextValue "E10" "  -0.24%  "
Set-TextValue "D11" "1.056"
Set-TextValue "E11" "  +1.49%  "
Set-TextValue "D12" "1.003"
Set-TextValue "E12" "  +0.21%  "
Set-TextValue "D13" "20.12"
Set-TextValue "E13" "  +1.86%  "
Set-TextValue "D14" "5.927"
Set-TextValue "E14" "  +1.80%  "
Set-TextValue "D15" "1.735.12"
Set-TextValue "E15" "  +0.75%  "
Set-TextValue "D16" "6.802"
Set-TextValue "E16" "  -0.92%  "
Set-TextValue "D17" "86.31"
Set-TextValue "E17" "  -2.55%  "
Set-TextValue "D18" "0.00001032"
Set-TextValue "E18" "  -0.67%  "
Set-TextValue "D19" "0.06417"
Set-TextValue "E19" "  +1.40%  "
Set-TextValue "D20" "1.002"
Set-TextValue "E20" "  +0.03%  "
Set-TextValue "D21" "16.42"
Set-TextValue "E21" "  -0.42%  "
Set-TextValue "D22" "5.718"
Set-TextValue "E22" "  +1.96%  "
Set-TextValue "D23" "26.939.78"
Set-TextValue "E23" "  -0.79%  "
Set-TextValue "D24" "11.22"
Set-TextValue "E24" "  +3.34%  "
Set-TextValue "D25" "2.050"
Set-TextValue "E25" "  -3.58%  "
Set-TextValue "D26" "153.80"
Set-TextValue "E26" "  -0.69%  "
Set-TextValue "D27" "19.83"
Set-TextValue "E27" "  +2.79%  "
Set-TextValue "D28" "1.935.49"
Set-TextValue "E28" "  +1.00%  "
Set-TextValue "D29" "2.201"
Set-TextValue "E29" "  +3.58%  "
Set-TextValue "D30" "119.70"
Set-TextValue "E30" "  -0.03%  "
Set-TextValue "D31" "1.041"
Set-TextValue "E31" "  +1.74%  "
Set-TextValue "D32" "0.09505"
Set-TextValue "E32" "  +4.54%  "
Set-TextValue "D33" "3.580"
Set-TextValue "E33" "  -0.51%  "
Set-TextValue "D34" "5.345"
Set-TextValue "E34" "  +0.14%  "
Set-TextValue "D35" "0.02188"
Set-TextValue "E35" "  -0.29%  "
Set-TextValue "D36" "0.05832"
Set-TextValue "E36" "  -0.29%  "
Set-TextValue "D37" "11.01"
Set-TextValue "E37" "  -0.22%  "
Set-TextValue "D38" "0.1996"
Set-TextValue "E38" "  +0.15%  "
Set-TextValue "B39" "WEMIXTOKEN"
Set-TextValue "C39" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D39" "1.420"
Set-TextValue "E39" "  +1.57%  "
Set-TextValue "B40" "InternetComputer(DFINITY)"
Set-TextValue "C40" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D40" "4.760"
Set-TextValue "E40" "  +1.03%  "
Set-TextValue "D41" "0.6024"
Set-TextValue "E41" "  +1.90%  "
Set-TextValue "D42" "1.106"
Set-TextValue "E42" "  -1.89%  "
Set-TextValue "D43" "7.594"
Set-TextValue "E43" "  +1.84%  "
Set-TextValue "D44" "12.81"
Set-TextValue "E44" "  +1.00%  "
Set-TextValue "D45" "3.590"
Set-TextValue "E45" "  +0.12%  "
Set-TextValue "D46" "0.5627"
Set-TextValue "E46" "  -0.13%  "
Set-TextValue "D47" "119.71"
Set-TextValue "E47" "  +0.74%  "
Set-TextValue "D48" "1.842"
Set-TextValue "E48" "  -1.06%  "
Set-TextValue "D49" "0.06657"
Set-TextValue "E49" "  +0.08%  "
Set-TextValue "D50" "1.097"
Set-TextValue "E50" "  +1.59%  "
Set-TextValue "E51" "  -0.21%  "
